$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate alt_label value in C11 (landings_expected duplicate)
$ws.Range("C11").ClearContents()

# Rename log_R0 -> log_recruitment_unfished (used as alt_label for log_rzero, row 30)
$ws.Range("C30").Value = "log_recruitment_unfished"

# Correct alt_label for IndexNumberAtAge / IndexNumberAtLength: indices -> indices_numbers
$ws.Range("C38").Value = "indices_numbers"
$ws.Range("C39").Value = "indices_numbers"

# Rename q -> catchability (used as alt_label for Q, row 46)
$ws.Range("C46").Value = "catchability"

# Update the view: scroll position, zoom, and selection
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("C40").Select()
